# Insert a new data row (shift existing rows 120..231 down to 121..232)
# and populate the new row 120 with the latest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 120:231 down by one row, creating a blank row 120.
$ws.Rows("120:120").Insert()

# Fill the newly inserted row 120 with the new record.
$ws.Range("A120").Value2 = 3
$ws.Range("B120").Value2 = "Femacal de La Calera"
$ws.Range("C120").Value2 = "Coquimbo"
$ws.Range("D120").Value2 = 44512
$ws.Range("E120").Value2 = 5
$ws.Range("F120").Value2 = 100112043
$ws.Range("G120").Value2 = "Pepino ensalada"
$ws.Range("H120").Value2 = "Sin especificar"
$ws.Range("I120").Value2 = "Primera"
$ws.Range("J120").Value2 = 110
$ws.Range("K120").Value2 = 7500
$ws.Range("L120").Value2 = 8000
$ws.Range("M120").Value2 = 7727
$ws.Range("N120").Value2 = "`$/caja 70 unidades"
$ws.Range("O120").Value2 = "Región de Arica y Parinacota"
$ws.Range("P120").Value2 = 110
$ws.Range("Q120").Value2 = 70
$ws.Range("R120").Value2 = "Hortaliza"
